$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header style from an existing header cell (e.g. E1) into F1:H1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Boolean data, default FALSE for all rows 2-12, columns F:H
$values = @{
    2  = @($false, $false, $false)
    3  = @($false, $false, $false)
    4  = @($false, $false, $false)
    5  = @($false, $true,  $false)
    6  = @($false, $false, $false)
    7  = @($false, $false, $false)
    8  = @($false, $false, $false)
    9  = @($false, $false, $false)
    10 = @($false, $false, $false)
    11 = @($false, $true,  $false)
    12 = @($false, $false, $false)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("F$row").Value = $rowVals[0]
    $ws.Range("G$row").Value = $rowVals[1]
    $ws.Range("H$row").Value = $rowVals[2]
}
